$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "mobility" data-description rows appended after the existing data
# (rows 59-64), each following the same layout as the rest of the sheet:
#   A = Header/Field Name     B = Format (csv)
#   C = Category (Environmental Data)
#   D = Description           E = Context Area (Indonesia)
$rows = @(
    @{ Row = 59; Name = "Retail and Recreation Mobility"; Desc = "Percent change from baseline for mobility in retail and recreation" },
    @{ Row = 60; Name = "Grocery and Pharmacy Mobility";  Desc = "Percent change from baseline for mobility in Grocery and Pharmacy" },
    @{ Row = 61; Name = "Parks Mobility";                 Desc = "Percent change from baseline for mobility in parks" },
    @{ Row = 62; Name = "Transit Mobility";                Desc = "Percent change from baseline for mobility for transit" },
    @{ Row = 63; Name = "Workplace Mobility";              Desc = "Percent change from baseline for mobility in workplaces" },
    @{ Row = 64; Name = "Residential Mobility";            Desc = "Percent change from baseline for mobility in residences" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Name
    $ws.Cells.Item($rowNum, 2).Value = "csv"
    $ws.Cells.Item($rowNum, 3).Value = "Environmental Data"
    $ws.Cells.Item($rowNum, 4).Value = $r.Desc
    $ws.Cells.Item($rowNum, 4).WrapText = $true
    $ws.Cells.Item($rowNum, 5).Value = "Indonesia"
}

# The two longest descriptions wrap to two lines at the existing column D
# width, same as rows 57/58.
$ws.Rows(59).RowHeight = 25.5
$ws.Rows(60).RowHeight = 25.5

# Column A needs to widen to fit the new, longer header names.
$ws.Columns("A").ColumnWidth = 32.14

# Leave the cursor on the last entered cell, matching where editing ended.
$ws.Range("D64").Select()
